# Update scripts with new TPM values (Dll4-Notch1.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ E=3; F=1; G=36.75793933333333; H=110.273818; I=0.9858943139827973; J=0.9858943139827971; M=38.55267666666666; N=115.65803; O=0.5758151725879548; P=0.5758151725879548; Q=1417.116950050949; R=12754.05255045854; S=0.5676929045594877; T=0.5676929045594876 }
    3  = @{ E=3; F=1; G=36.75793933333333; H=110.273818; I=0.9858943139827973; J=0.9858943139827971; O=0.08021535714867321; P=0.08021535714867323; Q=197.4149825869893; R=1776.734843282904; S=0.07908386450697626; T=0.07908386450697626 }
    4  = @{ E=3; F=1; G=36.75793933333333; H=110.273818; I=0.9858943139827973; J=0.9858943139827971; M=23.02986166666667; N=69.089585; O=0.3439694702633719; P=0.3439694702633719; Q=846.5302579983921; R=7618.772321985531; S=0.3391175449163332; T=0.3391175449163332 }
    5  = @{ I=0.001251989679428792; J=0.001251989679428792; M=38.55267666666666; N=115.65803; O=0.5758151725879548; P=0.5758151725879548; Q=1.799600394123333; R=16.19640354711; S=0.0007209146533386282; T=0.000720914653338628 }
    6  = @{ I=0.001251989679428792; J=0.001251989679428792; O=0.08021535714867321; P=0.08021535714867323; S=0.0001004287992818335; T=0.0001004287992818334 }
    7  = @{ I=0.001251989679428792; J=0.001251989679428792; M=23.02986166666667; N=69.089585; O=0.3439694702633719; P=0.3439694702633719; Q=1.075010912738333; R=9.675098214644999; S=0.0004306462268083305; T=0.0004306462268083304 }
    8  = @{ G=0.4792353333333333; H=1.437706; I=0.01285369633777395; J=0.01285369633777395; M=38.55267666666666; N=115.65803; O=0.5758151725879548; P=0.5758151725879548; Q=18.47580485324222; R=166.28224367918; S=0.007401353375128471; T=0.007401353375128469 }
    9  = @{ G=0.4792353333333333; H=1.437706; I=0.01285369633777395; J=0.01285369633777395; O=0.08021535714867321; P=0.08021535714867323; Q=2.573817703085334; R=23.164359327768; S=0.001031063842415131; T=0.001031063842415131 }
    10 = @{ G=0.4792353333333333; H=1.437706; I=0.01285369633777395; J=0.01285369633777395; M=23.02986166666667; N=69.089585; O=0.3439694702633719; P=0.3439694702633719; Q=11.03672343244555; R=99.33051089201; S=0.004421279120230349; T=0.004421279120230349 }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
